$d = $word.ActiveDocument

# Replace the full contents of a paragraph (identified by 1-based index) with
# the supplied run-level WordprocessingML, preserving whatever runs (including
# any pre-existing empty <w:r/> marker runs) were not part of the edit.
function Set-ParagraphRuns {
    param(
        [int]$ParaIndex,
        [string]$InnerXml
    )
    $p = $d.Paragraphs.Item($ParaIndex)
    $r = $p.Range   # full paragraph range, including its paragraph mark
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p>' + $InnerXml + '</w:p></w:body>' +
           '</w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'
    [void]$r.InsertXML($pkg)
}

# Same idea but replaces two consecutive paragraphs in one shot. Needed for
# the trailing paragraphs of the document: InsertXML over a range that ends
# exactly at the document's final paragraph mark (the body's last <w:p>)
# otherwise leaves behind a stray trailing empty paragraph, so the last two
# paragraphs must be rewritten together.
function Set-TwoParagraphRuns {
    param(
        [int]$FirstParaIndex,
        [string]$FirstInnerXml,
        [int]$SecondParaIndex,
        [string]$SecondInnerXml
    )
    $p1 = $d.Paragraphs.Item($FirstParaIndex)
    $p2 = $d.Paragraphs.Item($SecondParaIndex)
    $r = $d.Range($p1.Range.Start, $p2.Range.End)
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p>' + $FirstInnerXml + '</w:p><w:p>' + $SecondInnerXml + '</w:p></w:body>' +
           '</w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'
    [void]$r.InsertXML($pkg)
}

# 1: Heading1 "Play ... for Free!" -> "...Free" (single run, no empty-run prefix)
Set-ParagraphRuns 1 '<w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Play Little Green Men Nova Wilds for Free</w:t></w:r>'

# 2: What we like - bullet 1
Set-ParagraphRuns 37 '<w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r/><w:r><w:t>Exciting and varied gameplay mechanics</w:t></w:r>'

# 3: What we like - bullet 2
Set-ParagraphRuns 38 '<w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r/><w:r><w:t>Cartoonish graphics style with bright colors</w:t></w:r>'

# 4: What we like - bullet 3
Set-ParagraphRuns 39 '<w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r/><w:r><w:t>Alien theme with unique symbols</w:t></w:r>'

# 5: What we like - bullet 4
Set-ParagraphRuns 40 '<w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r/><w:r><w:t>Multiple special functions and free spins</w:t></w:r>'

# 6: What we don't like - bullet 1
Set-ParagraphRuns 42 '<w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r/><w:r><w:t>RTP is average compared to other slots</w:t></w:r>'

# 7: What we don't like - bullet 2
Set-ParagraphRuns 43 '<w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r/><w:r><w:t>Limited number of paylines</w:t></w:r>'

# 8 & 9: Bold CTA + Italic summary line (last two paragraphs in the body -
# must be rewritten together, see Set-TwoParagraphRuns comment above)
Set-TwoParagraphRuns 44 '<w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Little Green Men Nova Wilds for Free</w:t></w:r>' `
                     45 '<w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Read our review of Little Green Men Nova Wilds and play for free to experience the exciting gameplay and unique alien theme.</w:t></w:r>'
